$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: "Both Side" command renamed to "Drive" ---
$ws.Range("B9").Value = "Drive"

# --- Row 10: new "Steer" command row ---
$ws.Range("B10").Value = "Steer"
$ws.Range("C10").Value = 9002
$ws.Range("D10").Value = "S"
$ws.Range("E10").Value = "N/A"
$ws.Range("G10").Value = "Range is ±90"
$ws.Range("F10").Value = "steer left, steer right"
$ws.Range("H10").Value = """S45,-45"""

# --- Row 11: new "Steer Feedback" command row ---
$ws.Range("B11").Value = "Steer Feedback"
$ws.Range("C11").Value = 9004
$ws.Range("D11").Value = "N/A"
$ws.Range("E11").Value = "N/A"
$ws.Range("F11").Value = "N/A"
$ws.Range("G11").Value = "Simply open the connection to receive data"
$ws.Range("G11").Style = "Normal"
$ws.Range("I11").Value = "left front, left back, right front, right back"
$ws.Range("J11").Value = "0.0,0.0,0.0,0.0"

# --- Rows 30-33: fill in Port numbers for the Autonomous commands ---
$ws.Range("C30").Value = 9004
$ws.Range("C31").Value = 9003
$ws.Range("C32").Value = 9003
$ws.Range("C33").Value = 9003

# --- Update selection / view to match the saved workbook state ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C35").Select()
